$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 269, shifting existing rows 269-284 down to 270-285.
$ws.Rows.Item(269).Insert()

# Copy the date-column formatting from the cell below (D270, the row that used to be D269)
# into the new D269 cell, so the date column keeps its date style.
$ws.Cells.Item(270, 4).Copy()
$ws.Cells.Item(269, 4).PasteSpecial(-4122) # xlPasteFormats

# Populate the new row 269 with the new weekly record.
$ws.Cells.Item(269, 1).Value = 8
$ws.Cells.Item(269, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(269, 3).Value = "Coquimbo"
$ws.Cells.Item(269, 4).Value = 45147
$ws.Cells.Item(269, 5).Value = 4
$ws.Cells.Item(269, 6).Value = 100112001
$ws.Cells.Item(269, 7).Value = "Berenjena"
$ws.Cells.Item(269, 8).Value = "Sin especificar"
$ws.Cells.Item(269, 9).Value = "Primera"
$ws.Cells.Item(269, 10).Value = 400
$ws.Cells.Item(269, 11).Value = 10000
$ws.Cells.Item(269, 12).Value = 11000
$ws.Cells.Item(269, 13).Value = 10500
$ws.Cells.Item(269, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(269, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(269, 16).Value = 210
$ws.Cells.Item(269, 17).Value = 50
$ws.Cells.Item(269, 18).Value = "Hortaliza"
